$d = $word.ActiveDocument

# 1. Title / H1 heading (also matches the bold "Play Ed Jones..." run near the
#    end of the document, which changes identically) -> replace all occurrences.
$d.Content.Find.Execute(
    "Play Ed Jones and Book of Bastet for Free - Slot Game Review",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Free Ed Jones and Book of Bastet Slot Game", 2)

# 2. "What we like" bullet list
$d.Content.Find.Execute(
    "Multiple bonus features",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Multiple bonus features for increased excitement", 2)

$d.Content.Find.Execute(
    "Intuitive and user-friendly interface",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Immersive and realistic graphics and sound", 2)

$d.Content.Find.Execute(
    "Attractive graphics and sound",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Wide range of betting options for all types of players", 2)

$d.Content.Find.Execute(
    "Variety of betting options",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Book of Bastet symbol serving as both Wild and Scatter", 2)

# 3. "What we don't like" bullet list.
#    NOTE: the two bullets "rotate" their text ("Gamble feature can be risky"
#    is both the NEW text of bullet 1 and the OLD text of bullet 2), so the
#    second bullet must be renamed away BEFORE the first bullet is given the
#    text "Gamble feature can be risky" - otherwise a later ReplaceAll would
#    catch both occurrences.
$d.Content.Find.Execute(
    "Gamble feature can be risky",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Limited number of paylines", 2)

$d.Content.Find.Execute(
    "High variance can lead to longer periods without a win",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Gamble feature can be risky", 2)

# 4. Closing italic summary paragraph.
$d.Content.Find.Execute(
    "Discover the ancient Egyptian world with Ed Jones and Book of Bastet. Play for free and read our review of this 5-reel, 10-payline slot game.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Ed Jones and Book of Bastet slot game and play for free.", 2)
